$wb = $excel.ActiveWorkbook

# "Overview" sheet: row 3 (7f1683b4-...) -> Latest HO Xliff Generate Date (col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(3, 7).Value = "2016-10-21 03:43:47"

# "zh-cn" sheet: row 3 (7f1683b4-...) -> Correspond Handoff Datetime (col H) / Correspond Handback DateTime (col K)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(3, 8).Value = "2016-10-21 03:43:36"
$wsZhCn.Cells.Item(3, 11).Value = "2016-10-21 03:44:19"

# "de-de" sheet: row 3 (7f1683b4-...) -> Correspond Handoff Datetime (col H, shared text w/ Overview) / Correspond Handback DateTime (col K)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(3, 8).Value = "2016-10-21 03:43:47"
$wsDeDe.Cells.Item(3, 11).Value = "2016-10-21 03:44:37"
